# Regenerate the handback-status report: new source-file GUIDs and new
# handoff/handback correlation file names + timestamps for the zh-cn and
# de-de rows, on the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# New values used everywhere below
# ---------------------------------------------------------------------
$guidRow2 = "bbfaa0a5-8a94-4e31-a525-b54c93ac7a79"
$guidRow3 = "ffff2536c0d5-fba3-4ed8-8872-a2f28643efd7"
$corrHash = "0d371a7d070655dc48d1e410b1c7fc550d25340f"

$mdRow2 = "$guidRow2.md"
$mdRow3 = "$guidRow3.md"

$zhXlf = "$guidRow2.$corrHash.zh-cn.xlf"
$deXlf = "$guidRow2.$corrHash.de-de.xlf"

$zhTime1 = "2016-03-21 17:06:23"
$zhTime2 = "2016-03-21 17:06:50"
$deTime1 = "2016-03-21 17:06:27"
$deTime2 = "2016-03-21 17:06:56"

# Hyperlink colour used through the workbook for the "HyperLink" look
# (matches the font already defined in styles.xml for style index 1).
$hlColor = 15570276   # BGR encoding of RGB FF6495ED

function Set-CellHyperlink($ws, $cellRef, $url, $display) {
    $rng = $ws.Range($cellRef)
    $ws.Hyperlinks.Add($rng, $url, "", "", $display) | Out-Null
    $rng.Font.Underline = 2
    $rng.Font.Color = $hlColor
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$ovUrlA2 = "https://github.com/OpenLocalizationTest/oltest/blob/aca6a5661060e972f397f826f0f2ad2665dce319/e2e/6107b258-e6c2-4976-89cc-cd2ddad365ee.md"
$ovUrlA3 = "https://github.com/OpenLocalizationTest/oltest/blob/aca6a5661060e972f397f826f0f2ad2665dce319/e2e/d588e240-9d59-4c20-8405-e3b8dc4c56a0.md"

$wsOverview.Hyperlinks.Delete() | Out-Null

Set-CellHyperlink $wsOverview "A2" $ovUrlA2 $mdRow2
Set-CellHyperlink $wsOverview "A3" $ovUrlA3 $mdRow3

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhUrlA2 = "https://github.com/OpenLocalizationTest/oltest/blob/aca6a5661060e972f397f826f0f2ad2665dce319/e2e/6107b258-e6c2-4976-89cc-cd2ddad365ee.md"
$zhUrlD2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/40567fd30c58f7b919472148ba1523fb6b56246d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/6107b258-e6c2-4976-89cc-cd2ddad365ee.af704c3e906db2eeb6280ec0bae964558a262c3a.zh-cn.xlf"
$zhUrlF2 = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/13fc9b7d247d26dc875e551bedbc487fc767b1af/e2e/6107b258-e6c2-4976-89cc-cd2ddad365ee.md"
$zhUrlG2 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e0a6ae18602cb4f31dc436f1e3445768cf6ca371/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/6107b258-e6c2-4976-89cc-cd2ddad365ee.af704c3e906db2eeb6280ec0bae964558a262c3a.zh-cn.xlf"
$zhUrlA3 = "https://github.com/OpenLocalizationTest/oltest/blob/aca6a5661060e972f397f826f0f2ad2665dce319/e2e/d588e240-9d59-4c20-8405-e3b8dc4c56a0.md"
$zhUrlD3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/40567fd30c58f7b919472148ba1523fb6b56246d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d588e240-9d59-4c20-8405-e3b8dc4c56a0.6194c2004d7ed5ea0a280caf6ac0255a29905edc.zh-cn.xlf"
$zhUrlF3 = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/13fc9b7d247d26dc875e551bedbc487fc767b1af/e2e/d588e240-9d59-4c20-8405-e3b8dc4c56a0.md"
$zhUrlG3 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e0a6ae18602cb4f31dc436f1e3445768cf6ca371/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d588e240-9d59-4c20-8405-e3b8dc4c56a0.6194c2004d7ed5ea0a280caf6ac0255a29905edc.zh-cn.xlf"

$wsZh.Hyperlinks.Delete() | Out-Null

Set-CellHyperlink $wsZh "A2" $zhUrlA2 $mdRow2
Set-CellHyperlink $wsZh "D2" $zhUrlD2 $zhXlf
Set-CellHyperlink $wsZh "F2" $zhUrlF2 $mdRow2
Set-CellHyperlink $wsZh "G2" $zhUrlG2 $zhXlf
$wsZh.Range("E2").Value = $zhTime1
$wsZh.Range("H2").Value = $zhTime2

Set-CellHyperlink $wsZh "A3" $zhUrlA3 $mdRow3
Set-CellHyperlink $wsZh "D3" $zhUrlD3 $zhXlf
Set-CellHyperlink $wsZh "F3" $zhUrlF3 $mdRow3
Set-CellHyperlink $wsZh "G3" $zhUrlG3 $zhXlf
$wsZh.Range("E3").Value = $zhTime1
$wsZh.Range("H3").Value = $zhTime2

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deUrlA2 = "https://github.com/OpenLocalizationTest/oltest/blob/aca6a5661060e972f397f826f0f2ad2665dce319/e2e/6107b258-e6c2-4976-89cc-cd2ddad365ee.md"
$deUrlD2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a9d28c8fe92b1642bd09f50827b0a4588d177304/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/6107b258-e6c2-4976-89cc-cd2ddad365ee.af704c3e906db2eeb6280ec0bae964558a262c3a.de-de.xlf"
$deUrlF2 = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/a1316e7715a2f73829f24fc9d4f1f3f777a42f8d/e2e/6107b258-e6c2-4976-89cc-cd2ddad365ee.md"
$deUrlG2 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c39528c5bb7c1ae86658be2e3cdd8d838559239a/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/6107b258-e6c2-4976-89cc-cd2ddad365ee.af704c3e906db2eeb6280ec0bae964558a262c3a.de-de.xlf"
$deUrlA3 = "https://github.com/OpenLocalizationTest/oltest/blob/aca6a5661060e972f397f826f0f2ad2665dce319/e2e/d588e240-9d59-4c20-8405-e3b8dc4c56a0.md"
$deUrlD3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a9d28c8fe92b1642bd09f50827b0a4588d177304/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d588e240-9d59-4c20-8405-e3b8dc4c56a0.6194c2004d7ed5ea0a280caf6ac0255a29905edc.de-de.xlf"
$deUrlF3 = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/a1316e7715a2f73829f24fc9d4f1f3f777a42f8d/e2e/d588e240-9d59-4c20-8405-e3b8dc4c56a0.md"
$deUrlG3 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c39528c5bb7c1ae86658be2e3cdd8d838559239a/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d588e240-9d59-4c20-8405-e3b8dc4c56a0.6194c2004d7ed5ea0a280caf6ac0255a29905edc.de-de.xlf"

$wsDe.Hyperlinks.Delete() | Out-Null

Set-CellHyperlink $wsDe "A2" $deUrlA2 $mdRow2
Set-CellHyperlink $wsDe "D2" $deUrlD2 $deXlf
Set-CellHyperlink $wsDe "F2" $deUrlF2 $mdRow2
Set-CellHyperlink $wsDe "G2" $deUrlG2 $deXlf
$wsDe.Range("E2").Value = $deTime1
$wsDe.Range("H2").Value = $deTime2

Set-CellHyperlink $wsDe "A3" $deUrlA3 $mdRow3
Set-CellHyperlink $wsDe "D3" $deUrlD3 $deXlf
Set-CellHyperlink $wsDe "F3" $deUrlF3 $mdRow3
Set-CellHyperlink $wsDe "G3" $deUrlG3 $deXlf
$wsDe.Range("E3").Value = $deTime1
$wsDe.Range("H3").Value = $deTime2

Write-Host "Handback status report regenerated."
